$wb = $excel.ActiveWorkbook

# --- Sheet "Valeurs réelles": update header labels (add _class suffix) and classifier outputs ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value = "PRIX EXP POMME FUJI FRANCE 201/270G CAT.I PLATEAU 1RG_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME FUJI FRANCE 201/270G CAT.I PLATEAU 1RG_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME FUJI FRANCE 201/270G CAT.I PLATEAU 1RG_S+3_class"

$sheet1Data = @{
    2 = @(1, 2, 2)
    3 = @(2, 2, 2)
    4 = @(2, 2, 2)
    5 = @(2, 2, 2)
    6 = @(2, 2, 2)
    7 = @(2, 2, 2)
    8 = @(2, 2, 2)
    9 = @(2, 2, 2)
    10 = @(2, 2, 2)
    11 = @(2, 2, 2)
    12 = @(2, 2, 2)
    13 = @(2, 2, 2)
    14 = @(2, 2, 2)
    15 = @(2, 2, 2)
    16 = @(2, 2, 2)
    17 = @(2, 2, 2)
    18 = @(2, 2, 2)
    19 = @(2, 2, 2)
    20 = @(2, 2, 2)
    21 = @(2, 2, 2)
    22 = @(2, 2, 1)
    23 = @(2, 1, 0)
    24 = @(1, 0, 4)
    25 = @(0, 4, 2)
    26 = @(4, 2, 2)
    27 = @(2, 2, 2)
    28 = @(2, 2, 2)
}

foreach ($r in $sheet1Data.Keys) {
    $vals = $sheet1Data[$r]
    $ws1.Cells.Item($r, 3).Value = $vals[0]
    $ws1.Cells.Item($r, 4).Value = $vals[1]
    $ws1.Cells.Item($r, 5).Value = $vals[2]
}

# --- Sheet "Prédictions": replace regressor outputs with classifier outputs ---
$ws2 = $wb.Worksheets.Item(2)

$sheet2Data = @{
    2 = @(0, 0, 0)
    3 = @(0, 0, 0)
    4 = @(0, 0, 0)
    5 = @(0, 0, 0)
    6 = @(0, 0, 0)
    7 = @(0, 0, 0)
    8 = @(0, 0, 0)
    9 = @(0, 0, 0)
    10 = @(0, 0, 2)
    11 = @(0, 0, 0)
    12 = @(0, 0, 0)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(-2, 2, 0)
    16 = @(-2, 2, -1)
    17 = @(0, -1, 0)
    18 = @(0, 2, 0)
    19 = @(-1, -1, 0)
    20 = @(0, -1, 0)
    21 = @(-1, -1, 0)
    22 = @(0, -1, 0)
    23 = @(0, 0, 0)
    24 = @(0, 0, 0)
    25 = @(0, 0, 1)
    26 = @(0, 0, 0)
    27 = @(0, 0, 0)
    28 = @(0, -1, 0)
}

foreach ($r in $sheet2Data.Keys) {
    $vals = $sheet2Data[$r]
    $ws2.Cells.Item($r, 2).Value = $vals[0]
    $ws2.Cells.Item($r, 3).Value = $vals[1]
    $ws2.Cells.Item($r, 4).Value = $vals[2]
}
